$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "master 2 commit"
$ws.Range("A2").Select()
